$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data table (rows 2-7), following Dr Hou's advice:
# Sending cluster x Target cluster combinations for Efnb3 -> Ephb2
$data = @(
    @("ECs",  "Efnb3", "Ephb2", "FAPs", 1, 0.3333333333333333, 0.09174599999999999, 0.275238, 0.08724331438250911, 0.08724331438250911, 3, 1, 6.346253666666667, 19.038761, 0.9446330608455225, 0.9446330608455226, 0.582243388902, 5.240190500118, 0.08241291910345777, 0.08241291910345779),
    @("ECs",  "Efnb3", "Ephb2", "sCs",  1, 0.3333333333333333, 0.09174599999999999, 0.275238, 0.08724331438250911, 0.08724331438250911, 3, 1, 0.3719673333333333, 1.115902, 0.05536693915447755, 0.05536693915447755, 0.03412651496399999, 0.3071386346759999, 0.004830395279051338, 0.004830395279051338),
    @("FAPs", "Efnb3", "Ephb2", "FAPs", 2, 0.6666666666666666, 0.1498043333333333, 0.449413, 0.1424522763811195, 0.1424522763811195, 3, 1, 6.346253666666667, 19.038761, 0.9446330608455225, 0.9446330608455226, 0.9506962996992223, 8.556266697293001, 0.1345651298623093, 0.1345651298623093),
    @("FAPs", "Efnb3", "Ephb2", "sCs",  2, 0.6666666666666666, 0.1498043333333333, 0.449413, 0.1424522763811195, 0.1424522763811195, 3, 1, 0.3719673333333333, 1.115902, 0.05536693915447755, 0.05536693915447755, 0.05572231839177778, 0.501500865526, 0.007887146518810263, 0.007887146518810263),
    @("sCs",  "Efnb3", "Ephb2", "FAPs", 3, 1, 0.8100603333333333, 2.430181, 0.7703044092363713, 0.7703044092363713, 3, 1, 6.346253666666667, 19.038761, 0.9446330608455225, 0.9446330608455226, 5.140848360637889, 46.267635245741, 0.7276550118797553, 0.7276550118797555),
    @("sCs",  "Efnb3", "Ephb2", "sCs",  3, 1, 0.8100603333333333, 2.430181, 0.7703044092363713, 0.7703044092363713, 3, 1, 0.3719673333333333, 1.115902, 0.05536693915447755, 0.05536693915447755, 0.3013159820291111, 2.711843838261999, 0.04264939735661595, 0.04264939735661595)
)

$rowIndex = 2
foreach ($row in $data) {
    $colIndex = 1
    foreach ($value in $row) {
        $ws.Cells.Item($rowIndex, $colIndex).Value = $value
        $colIndex++
    }
    $rowIndex++
}
